$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E and F hold numeric-looking text (IDs), keep them as text.
$ws.Range("E2:F3").NumberFormat = "@"

# Update existing row 2 values
$ws.Range("B2").Value = "уоклшупк"
$ws.Range("C2").Value = "Иванов"
$ws.Range("D2").Value = "Иван"
$ws.Range("E2").Value = "4315"
$ws.Range("F2").Value = "1008919333"

# Add new row 3
$ws.Range("A3").Value = 7918574921
$ws.Range("B3").Value = "халоо"
$ws.Range("C3").Value = "Иванов"
$ws.Range("D3").Value = "Иван"
$ws.Range("E3").Value = "4319"
$ws.Range("F3").Value = "1008919333"

# Update selection to match target (A2 active cell, A2:F3 selected range)
$ws.Range("A2:F3").Select()
